$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf17"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1403416666666667
$ws.Range("H2").Value = 0.421025
$ws.Range("I2").Value = 0.1831816192293602
$ws.Range("J2").Value = 0.1831816192293602
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4581623333333333
$ws.Range("N2").Value = 1.374487
$ws.Range("O2").Value = 0.1060599910922654
$ws.Range("P2").Value = 0.1060599910922654
$ws.Range("Q2").Value = 0.0642992654638889
$ws.Range("R2").Value = 0.578693389175
$ws.Range("S2").Value = 0.0194282409037327
$ws.Range("T2").Value = 0.0194282409037327
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf17"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1403416666666667
$ws.Range("H3").Value = 0.421025
$ws.Range("I3").Value = 0.1831816192293602
$ws.Range("J3").Value = 0.1831816192293602
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.075101
$ws.Range("N3").Value = 0.225303
$ws.Range("O3").Value = 0.01738512926863672
$ws.Range("P3").Value = 0.01738512926863672
$ws.Range("Q3").Value = 0.01053979950833333
$ws.Range("R3").Value = 0.094858195575
$ws.Range("S3").Value = 0.003184636129940617
$ws.Range("T3").Value = 0.003184636129940617
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf17"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1403416666666667
$ws.Range("H4").Value = 0.421025
$ws.Range("I4").Value = 0.1831816192293602
$ws.Range("J4").Value = 0.1831816192293602
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.786578
$ws.Range("N4").Value = 11.359734
$ws.Range("O4").Value = 0.8765548796390978
$ws.Range("P4").Value = 0.8765548796390978
$ws.Range("Q4").Value = 0.5314146674833333
$ws.Range("R4").Value = 4.78273200735
$ws.Range("S4").Value = 0.1605687421956869
$ws.Range("T4").Value = 0.1605687421956869
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf17"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.366566
$ws.Range("H5").Value = 1.099698
$ws.Range("I5").Value = 0.4784619922885553
$ws.Range("J5").Value = 0.4784619922885552
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4581623333333333
$ws.Range("N5").Value = 1.374487
$ws.Range("O5").Value = 0.1060599910922654
$ws.Range("P5").Value = 0.1060599910922654
$ws.Range("Q5").Value = 0.1679467338806667
$ws.Range("R5").Value = 1.511520604926
$ws.Range("S5").Value = 0.05074567464011174
$ws.Range("T5").Value = 0.05074567464011174
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf17"
$ws.Range("C6").Value = "Fgfr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.366566
$ws.Range("H6").Value = 1.099698
$ws.Range("I6").Value = 0.4784619922885553
$ws.Range("J6").Value = 0.4784619922885552
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.075101
$ws.Range("N6").Value = 0.225303
$ws.Range("O6").Value = 0.01738512926863672
$ws.Range("P6").Value = 0.01738512926863672
$ws.Range("Q6").Value = 0.027529473166
$ws.Range("R6").Value = 0.247765258494
$ws.Range("S6").Value = 0.008318123586065998
$ws.Range("T6").Value = 0.008318123586065998
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf17"
$ws.Range("C7").Value = "Fgfr4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.366566
$ws.Range("H7").Value = 1.099698
$ws.Range("I7").Value = 0.4784619922885553
$ws.Range("J7").Value = 0.4784619922885552
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.786578
$ws.Range("N7").Value = 11.359734
$ws.Range("O7").Value = 0.8765548796390978
$ws.Range("P7").Value = 0.8765548796390978
$ws.Range("Q7").Value = 1.388030751148
$ws.Range("R7").Value = 12.492276760332
$ws.Range("S7").Value = 0.4193981940623775
$ws.Range("T7").Value = 0.4193981940623775
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf17"
$ws.Range("C8").Value = "Fgfr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.12712
$ws.Range("H8").Value = 0.38136
$ws.Range("I8").Value = 0.1659239767455824
$ws.Range("J8").Value = 0.1659239767455823
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.4581623333333333
$ws.Range("N8").Value = 1.374487
$ws.Range("O8").Value = 0.1060599910922654
$ws.Range("P8").Value = 0.1060599910922654
$ws.Range("Q8").Value = 0.05824159581333334
$ws.Range("R8").Value = 0.52417436232
$ws.Range("S8").Value = 0.01759789549562972
$ws.Range("T8").Value = 0.01759789549562972
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf17"
$ws.Range("C9").Value = "Fgfr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.12712
$ws.Range("H9").Value = 0.38136
$ws.Range("I9").Value = 0.1659239767455824
$ws.Range("J9").Value = 0.1659239767455823
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.075101
$ws.Range("N9").Value = 0.225303
$ws.Range("O9").Value = 0.01738512926863672
$ws.Range("P9").Value = 0.01738512926863672
$ws.Range("Q9").Value = 0.009546839120000001
$ws.Range("R9").Value = 0.08592155208
$ws.Range("S9").Value = 0.002884609784488222
$ws.Range("T9").Value = 0.002884609784488222
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf17"
$ws.Range("C10").Value = "Fgfr4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.12712
$ws.Range("H10").Value = 0.38136
$ws.Range("I10").Value = 0.1659239767455824
$ws.Range("J10").Value = 0.1659239767455823
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.786578
$ws.Range("N10").Value = 11.359734
$ws.Range("O10").Value = 0.8765548796390978
$ws.Range("P10").Value = 0.8765548796390978
$ws.Range("Q10").Value = 0.4813497953600001
$ws.Range("R10").Value = 4.33214815824
$ws.Range("S10").Value = 0.1454414714654644
$ws.Range("T10").Value = 0.1454414714654644
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Fgf17"
$ws.Range("C11").Value = "Fgfr4"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.1321063333333333
$ws.Range("H11").Value = 0.396319
$ws.Range("I11").Value = 0.1724324117365021
$ws.Range("J11").Value = 0.1724324117365021
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.4581623333333333
$ws.Range("N11").Value = 1.374487
$ws.Range("O11").Value = 0.1060599910922654
$ws.Range("P11").Value = 0.1060599910922654
$ws.Range("Q11").Value = 0.06052614592811111
$ws.Range("R11").Value = 0.544735313353
$ws.Range("S11").Value = 0.01828818005279126
$ws.Range("T11").Value = 0.01828818005279126
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Fgf17"
$ws.Range("C12").Value = "Fgfr4"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.1321063333333333
$ws.Range("H12").Value = 0.396319
$ws.Range("I12").Value = 0.1724324117365021
$ws.Range("J12").Value = 0.1724324117365021
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.075101
$ws.Range("N12").Value = 0.225303
$ws.Range("O12").Value = 0.01738512926863672
$ws.Range("P12").Value = 0.01738512926863672
$ws.Range("Q12").Value = 0.009921317739666665
$ws.Range("R12").Value = 0.089291859657
$ws.Range("S12").Value = 0.002997759768141881
$ws.Range("T12").Value = 0.002997759768141881
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Fgf17"
$ws.Range("C13").Value = "Fgfr4"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.1321063333333333
$ws.Range("H13").Value = 0.396319
$ws.Range("I13").Value = 0.1724324117365021
$ws.Range("J13").Value = 0.1724324117365021
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.786578
$ws.Range("N13").Value = 11.359734
$ws.Range("O13").Value = 0.8765548796390978
$ws.Range("P13").Value = 0.8765548796390978
$ws.Range("Q13").Value = 0.5002309354606667
$ws.Range("R13").Value = 4.502078419146
$ws.Range("S13").Value = 0.151146471915569
$ws.Range("T13").Value = 0.151146471915569
